$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.22"
$ws.Range("E2").Value = "'0.02%"
$ws.Range("D3").Value = "'36.19"
$ws.Range("E3").Value = "'-2.04%"
$ws.Range("D4").Value = "'5.031"
$ws.Range("E4").Value = "'-0.16%"
$ws.Range("D5").Value = "'0.07866"
$ws.Range("E5").Value = "'-0.14%"
$ws.Range("D6").Value = "'2.128"
$ws.Range("E6").Value = "'-3.59%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.141"
$ws.Range("E7").Value = "'2.59%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'7.954"
$ws.Range("E8").Value = "'-0.67%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9219"
$ws.Range("E9").Value = "'-0.70%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.09517"
$ws.Range("E10").Value = "'-3.68%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1849"
$ws.Range("E11").Value = "'-1.76%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08572"
$ws.Range("E12").Value = "'-1.44%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03607"
$ws.Range("E13").Value = "'-0.01%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09928"
$ws.Range("E14").Value = "'-0.25%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001427"
$ws.Range("E15").Value = "'-3.92%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005699"
$ws.Range("E16").Value = "'0.64%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.470"
$ws.Range("E17").Value = "'0.35%"
$ws.Range("D18").Value = "'2.753"
$ws.Range("E18").Value = "'10.88%"
$ws.Range("D19").Value = "'0.3373"
$ws.Range("E19").Value = "'-1.80%"
$ws.Range("D20").Value = "'0.1338"
$ws.Range("E20").Value = "'0.83%"
$ws.Range("D21").Value = "'5.165"
$ws.Range("E21").Value = "'7.94%"
$ws.Range("D22").Value = "'0.2250"
$ws.Range("E22").Value = "'2.29%"
$ws.Range("D23").Value = "'0.04582"
$ws.Range("E23").Value = "'-0.21%"
$ws.Range("D24").Value = "'0.001233"
$ws.Range("E24").Value = "'-1.40%"
$ws.Range("D25").Value = "'0.004795"
$ws.Range("E25").Value = "'-8.02%"
$ws.Range("E26").Value = "'-7.01%"
$ws.Range("D27").Value = "'0.0004753"
$ws.Range("E27").Value = "'74.88%"
$ws.Range("D39").Value = "'0.01853"
$ws.Range("E39").Value = "'1.41%"
$ws.Range("D40").Value = "'0.04703"
$ws.Range("E40").Value = "'-1.43%"
$ws.Range("D41").Value = "'0.007799"
$ws.Range("E41").Value = "'-0.41%"
$ws.Range("D42").Value = "'0.1387"
$ws.Range("E42").Value = "'-1.24%"
$ws.Range("D43").Value = "'0.007732"
$ws.Range("E43").Value = "'2.13%"
$ws.Range("D44").Value = "'0.002271"
$ws.Range("E44").Value = "'3.33%"
$ws.Range("D45").Value = "'0.01138"
$ws.Range("E45").Value = "'9.54%"
$ws.Range("D46").Value = "'0.00006356"
$ws.Range("E46").Value = "'0.26%"
$ws.Range("E47").Value = "'0.07%"
$ws.Range("E48").Value = "'0.18%"
$ws.Range("D49").Value = "'52.20"
$ws.Range("E49").Value = "'43.65%"
$ws.Range("D50").Value = "'0.001902"
$ws.Range("E50").Value = "'-29.30%"
$ws.Range("E51").Value = "'0.07%"
